$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the row above (row 54) down to the new row 55
$ws.Range("A54:H54").Copy()
$ws.Range("A55:H55").PasteSpecial(-4122)

$ws.Range("A55").Value = "Mohamed Raâfet"
$ws.Range("B55").Value = "Ben Khedher"
$ws.Range("C55").Value = "Université Laval"
$ws.Range("D55").Value = "Canada"
$ws.Range("E55").Value = "fd-CnyYAAAAJ"
$ws.Range("F55").Value = "M"
$ws.Range("G55").Value = 1987
$ws.Range("H55").Value = "Médecine, Biologie et Sciences de la Santé"

$ws.Range("A55:H55").Select()
